# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 119
$wsExhibit.Range("F4").Value = 1641
$wsExhibit.Range("F5").Value = 19
$wsExhibit.Range("F6").Value = 23
$wsExhibit.Range("F9").Value = 64

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 394
$wsAll.Range("F4").Value = 1641
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 419
$wsAll.Range("F8").Value = 0
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 512
